$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (after the weekly re-shuffle of entries), columns:
# D=Fecha, K=Variedad, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
# P=Precio promedio ponderado, R=Origen, S=Precio $/Kg

$ws.Range("D2").Value = 44301
$ws.Range("K2").Value = "Hachiya"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1139

$ws.Range("D4").Value = 44699
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 29000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1639

$ws.Range("D5").Value = 44313
$ws.Range("K5").Value = "Mankaki"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1194

$ws.Range("D6").Value = 45043
$ws.Range("K6").Value = "Fuyu"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 26000
$ws.Range("P6").Value = 25500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1417

$ws.Range("D7").Value = 44305
$ws.Range("K7").Value = "Mankaki"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1361

$ws.Range("D8").Value = 45071
$ws.Range("K8").Value = "Fuyu"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 110
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23455
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1303

$ws.Range("D9").Value = 44355
$ws.Range("K9").Value = "Mankaki"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 270
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1139
